$d = $word.ActiveDocument

function Replace-InScope($scopeRange, $old, $new) {
    $r = $scopeRange.Duplicate
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replacement failed for: $old"
    }
}

function InsertAfter-InScope($scopeRange, $anchor, $newLines) {
    $r = $scopeRange.Duplicate
    $ok = $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Anchor not found: $anchor"
    }
    $r.Collapse(0)
    foreach ($ln in $newLines) {
        $r.InsertAfter([char]11 + $ln)
        $r.Collapse(0)
    }
}

# Scope ranges for the two tables (paragraphs 2 and 4)
$table1 = $d.Paragraphs.Item(2).Range
$table2 = $d.Paragraphs.Item(4).Range

# ---- Table 1 (bkvlps) replacements ----
Replace-InScope $table1 "Dep. Variable:                 bkvlps   R-squared:                       0.028" "Dep. Variable:                 bkvlps   R-squared:                       0.006"
Replace-InScope $table1 "Model:                            OLS   Adj. R-squared:                  0.024" "Model:                            OLS   Adj. R-squared:                  0.005"
Replace-InScope $table1 "Method:                 Least Squares   F-statistic:                     6.779" "Method:                 Least Squares   F-statistic:                     9.927"
Replace-InScope $table1 "Date:                Mon, 27 May 2024   Prob (F-statistic):           1.54e-20" "Date:                Mon, 27 May 2024   Prob (F-statistic):           4.66e-37"
Replace-InScope $table1 "Time:                        00:40:22   Log-Likelihood:                -22202." "Time:                        03:02:12   Log-Likelihood:            -1.9917e+05"
Replace-InScope $table1 "No. Observations:                5158   AIC:                         4.445e+04" "No. Observations:               39266   AIC:                         3.984e+05"
Replace-InScope $table1 "Df Residuals:                    5135   BIC:                         4.460e+04" "Df Residuals:                   39241   BIC:                         3.986e+05"
Replace-InScope $table1 "Df Model:                          22                                         " "Df Model:                          24                                         "
Replace-InScope $table1 "const          12.1867      1.304      9.344      0.000       9.630      14.744" "const          12.5932      0.764     16.489      0.000      11.096      14.090"
Replace-InScope $table1 "motif_003    4.297e-15   5.93e-15      0.724      0.469   -7.33e-15    1.59e-14" "motif_003    1.912e-14   1.03e-14      1.854      0.064    -1.1e-15    3.93e-14"
Replace-InScope $table1 "motif_012   -7.658e-15   6.94e-15     -1.103      0.270   -2.13e-14    5.95e-15" "motif_012   -3.065e-14   1.87e-14     -1.635      0.102   -6.74e-14    6.09e-15"
Replace-InScope $table1 "motif_102     1.73e-15   5.76e-15      0.300      0.764   -9.56e-15     1.3e-14" "motif_102   -9.862e-15   1.25e-14     -0.788      0.431   -3.44e-14    1.47e-14"
Replace-InScope $table1 "motif_021D      0.0041      0.030      0.140      0.889      -0.054       0.062" "motif_021D     -0.0606      0.017     -3.615      0.000      -0.093      -0.028"
Replace-InScope $table1 "motif_021U     -1.3483      1.222     -1.103      0.270      -3.744       1.047" "motif_021U      0.1508      0.091      1.659      0.097      -0.027       0.329"
Replace-InScope $table1 "motif_021C     -0.4665      0.356     -1.310      0.190      -1.165       0.232" "motif_021C     -0.5259      0.321     -1.637      0.102      -1.156       0.104"
Replace-InScope $table1 "motif_111D     -0.0290      0.325     -0.089      0.929      -0.665       0.607" "motif_111D     -0.0964      0.030     -3.245      0.001      -0.155      -0.038"
Replace-InScope $table1 "motif_111U      0.4722      0.176      2.682      0.007       0.127       0.817" "motif_111U      0.3260      0.061      5.378      0.000       0.207       0.445"
Replace-InScope $table1 "motif_030T      5.3723      9.012      0.596      0.551     -12.294      23.039" "motif_030T  -3.532e-14   2.96e-14     -1.192      0.233   -9.34e-14    2.27e-14"
Replace-InScope $table1 "motif_030C   2.491e-15   1.74e-15      1.432      0.152    -9.2e-16     5.9e-15" "motif_030C  -5.384e-15    4.1e-15     -1.312      0.190   -1.34e-14    2.66e-15"
Replace-InScope $table1 "motif_201       0.7331      0.145      5.061      0.000       0.449       1.017" "motif_201       0.0466      0.009      5.337      0.000       0.030       0.064"
Replace-InScope $table1 "motif_120D      3.9097      7.350      0.532      0.595     -10.500      18.319" "motif_120D     -0.7944      0.619     -1.283      0.199      -2.008       0.419"
Replace-InScope $table1 "motif_120U      2.3404      2.180      1.074      0.283      -1.933       6.613" "motif_120U    -11.8350      6.486     -1.825      0.068     -24.548       0.878"
Replace-InScope $table1 "motif_120C     -3.6389      6.011     -0.605      0.545     -15.424       8.146" "motif_120C    -20.0111     15.823     -1.265      0.206     -51.025      11.002"
Replace-InScope $table1 "motif_210      16.4089      3.716      4.416      0.000       9.125      23.693" "motif_210      -4.7894      3.313     -1.446      0.148     -11.282       1.704"
Replace-InScope $table1 "motif_300      -3.2163      1.003     -3.206      0.001      -5.183      -1.249" "motif_300      -0.8635      0.213     -4.057      0.000      -1.281      -0.446"
Replace-InScope $table1 "year_2011.0    -1.0411      1.646     -0.633      0.527      -4.268       2.185" "year_2011.0     0.2969      1.073      0.277      0.782      -1.807       2.400"
Replace-InScope $table1 "year_2012.0    -0.6112      1.602     -0.382      0.703      -3.751       2.529" "year_2012.0     0.4267      1.070      0.399      0.690      -1.670       2.524"
Replace-InScope $table1 "year_2013.0    -0.3783      1.563     -0.242      0.809      -3.443       2.686" "year_2013.0     0.9911      1.063      0.933      0.351      -1.092       3.074"
Replace-InScope $table1 "year_2014.0    -0.5210      1.542     -0.338      0.735      -3.544       2.502" "year_2014.0     0.9449      1.052      0.898      0.369      -1.118       3.008"
Replace-InScope $table1 "year_2015.0    -0.0859      1.528     -0.056      0.955      -3.081       2.909" "year_2015.0     0.9495      1.051      0.903      0.366      -1.111       3.010"
Replace-InScope $table1 "year_2016.0     0.6859      1.526      0.450      0.653      -2.305       3.677" "year_2016.0     1.6593      1.056      1.571      0.116      -0.411       3.730"
Replace-InScope $table1 "year_2017.0     2.4555      1.523      1.612      0.107      -0.530       5.441" "year_2017.0     3.4483      1.057      3.264      0.001       1.377       5.519"
Replace-InScope $table1 "year_2018.0     3.3974      1.516      2.240      0.025       0.425       6.370" "year_2018.0     3.6095      1.055      3.421      0.001       1.542       5.677"
Replace-InScope $table1 "year_2019.0     3.9906      1.515      2.633      0.008       1.020       6.961" "year_2019.0     4.3761      1.051      4.162      0.000       2.315       6.437"
Replace-InScope $table1 "year_2020.0     5.0772      1.516      3.350      0.001       2.106       8.048" "year_2020.0     4.7005      1.043      4.507      0.000       2.656       6.745"
Replace-InScope $table1 "Omnibus:                     6585.118   Durbin-Watson:                   0.391" "Omnibus:                    75350.988   Durbin-Watson:                   0.417"
Replace-InScope $table1 "Prob(Omnibus):                  0.000   Jarque-Bera (JB):          2389542.500" "Prob(Omnibus):                  0.000   Jarque-Bera (JB):        579164556.938"
Replace-InScope $table1 "Skew:                           6.711   Prob(JB):                         0.00" "Skew:                          14.435   Prob(JB):                         0.00"
Replace-InScope $table1 "Kurtosis:                     107.586   Cond. No.                     1.25e+16" "Kurtosis:                     597.274   Cond. No.                     1.17e+16"
Replace-InScope $table1 "[2] The smallest eigenvalue is 4.33e-27. This might indicate that there are" "[2] The smallest eigenvalue is 6.47e-25. This might indicate that there are"

# ---- Table 1 (bkvlps) new year rows, inserted after year_2020.0 row ----
InsertAfter-InScope $table1 "year_2020.0     4.7005      1.043      4.507      0.000       2.656       6.745" @(
    "year_2021.0     4.9108      1.020      4.813      0.000       2.911       6.911",
    "year_2022.0     4.4013      1.022      4.307      0.000       2.398       6.404",
    "year_2023.0     7.1259      1.056      6.749      0.000       5.056       9.195"
)

# ---- Table 2 (epspx) replacements ----
Replace-InScope $table2 "Dep. Variable:                  epspx   R-squared:                       0.009" "Dep. Variable:                  epspx   R-squared:                       0.004"
Replace-InScope $table2 "Model:                            OLS   Adj. R-squared:                  0.004" "Model:                            OLS   Adj. R-squared:                  0.003"
Replace-InScope $table2 "Method:                 Least Squares   F-statistic:                     2.016" "Method:                 Least Squares   F-statistic:                     6.748"
Replace-InScope $table2 "Date:                Mon, 27 May 2024   Prob (F-statistic):            0.00329" "Date:                Mon, 27 May 2024   Prob (F-statistic):           2.18e-22"
Replace-InScope $table2 "Time:                        00:40:22   Log-Likelihood:                -14122." "Time:                        03:02:12   Log-Likelihood:            -1.3872e+05"
Replace-InScope $table2 "No. Observations:                5162   AIC:                         2.829e+04" "No. Observations:               39360   AIC:                         2.775e+05"
Replace-InScope $table2 "Df Residuals:                    5139   BIC:                         2.844e+04" "Df Residuals:                   39335   BIC:                         2.777e+05"
Replace-InScope $table2 "Df Model:                          22                                         " "Df Model:                          24                                         "
Replace-InScope $table2 "const           1.0425      0.272      3.837      0.000       0.510       1.575" "const           1.0501      0.162      6.468      0.000       0.732       1.368"
Replace-InScope $table2 "motif_003   -1.382e-16   1.44e-15     -0.096      0.924   -2.96e-15    2.69e-15" "motif_003    1.193e-15   6.26e-15      0.191      0.849   -1.11e-14    1.35e-14"
Replace-InScope $table2 "motif_012   -1.183e-16   1.34e-15     -0.088      0.930   -2.75e-15    2.52e-15" "motif_012   -2.839e-16   8.02e-16     -0.354      0.723   -1.86e-15    1.29e-15"
Replace-InScope $table2 "motif_102    2.022e-16   6.76e-16      0.299      0.765   -1.12e-15    1.53e-15" "motif_102    7.102e-16   3.46e-16      2.051      0.040    3.15e-17    1.39e-15"
Replace-InScope $table2 "motif_021D      0.0089      0.006      1.444      0.149      -0.003       0.021" "motif_021D     -0.0078      0.004     -2.189      0.029      -0.015      -0.001"
Replace-InScope $table2 "motif_021U     -0.3641      0.255     -1.431      0.153      -0.863       0.135" "motif_021U      0.0397      0.019      2.054      0.040       0.002       0.078"
Replace-InScope $table2 "motif_021C     -0.0552      0.074     -0.744      0.457      -0.201       0.090" "motif_021C     -0.0261      0.068     -0.382      0.703      -0.160       0.108"
Replace-InScope $table2 "motif_111D      0.0748      0.068      1.108      0.268      -0.058       0.207" "motif_111D     -0.0280      0.006     -4.431      0.000      -0.040      -0.016"
Replace-InScope $table2 "motif_111U     -0.0182      0.037     -0.497      0.620      -0.090       0.054" "motif_111U      0.0411      0.013      3.202      0.001       0.016       0.066"
Replace-InScope $table2 "motif_030T      0.5207      1.877      0.277      0.782      -3.160       4.201" "motif_030T   -2.17e-15   1.85e-15     -1.170      0.242   -5.81e-15    1.47e-15"
Replace-InScope $table2 "motif_030C   5.243e-16   1.64e-15      0.320      0.749   -2.69e-15    3.74e-15" "motif_030C  -1.006e-15   1.75e-15     -0.573      0.566   -4.44e-15    2.43e-15"
Replace-InScope $table2 "motif_201       0.1019      0.030      3.380      0.001       0.043       0.161" "motif_201       0.0098      0.002      5.289      0.000       0.006       0.013"
Replace-InScope $table2 "motif_120D      0.1778      1.531      0.116      0.908      -2.824       3.180" "motif_120D     -0.0150      0.131     -0.114      0.909      -0.273       0.243"
Replace-InScope $table2 "motif_120U      0.2528      0.454      0.557      0.578      -0.637       1.143" "motif_120U     -1.2996      1.379     -0.942      0.346      -4.003       1.404"
Replace-InScope $table2 "motif_120C      0.1147      1.252      0.092      0.927      -2.340       2.570" "motif_120C     -3.8458      3.365     -1.143      0.253     -10.442       2.750"
Replace-InScope $table2 "motif_210       0.1920      0.774      0.248      0.804      -1.325       1.709" "motif_210      -1.4938      0.705     -2.120      0.034      -2.875      -0.113"
Replace-InScope $table2 "motif_300      -0.5367      0.209     -2.568      0.010      -0.946      -0.127" "motif_300      -0.0997      0.045     -2.207      0.027      -0.188      -0.011"
Replace-InScope $table2 "year_2011.0    -0.0980      0.343     -0.286      0.775      -0.771       0.574" "year_2011.0     0.0097      0.228      0.042      0.966      -0.438       0.457"
Replace-InScope $table2 "year_2012.0    -0.0116      0.334     -0.035      0.972      -0.666       0.643" "year_2012.0     0.0695      0.227      0.306      0.760      -0.376       0.515"
Replace-InScope $table2 "year_2013.0    -0.0662      0.326     -0.203      0.839      -0.705       0.572" "year_2013.0     0.1516      0.226      0.671      0.502      -0.291       0.594"
Replace-InScope $table2 "year_2014.0    -0.0008      0.321     -0.003      0.998      -0.631       0.629" "year_2014.0     0.0221      0.224      0.099      0.921      -0.416       0.460"
Replace-InScope $table2 "year_2015.0    -0.2899      0.318     -0.911      0.362      -0.914       0.334" "year_2015.0    -0.4228      0.223     -1.893      0.058      -0.861       0.015"
Replace-InScope $table2 "year_2016.0    -0.0838      0.318     -0.264      0.792      -0.707       0.539" "year_2016.0    -0.1064      0.225     -0.474      0.636      -0.547       0.334"
Replace-InScope $table2 "year_2017.0     0.2955      0.317      0.932      0.352      -0.326       0.917" "year_2017.0     0.2096      0.225      0.933      0.351      -0.231       0.650"
Replace-InScope $table2 "year_2018.0     0.4360      0.316      1.380      0.168      -0.183       1.055" "year_2018.0     0.3484      0.224      1.553      0.120      -0.091       0.788"
Replace-InScope $table2 "year_2019.0     0.3533      0.316      1.119      0.263      -0.265       0.972" "year_2019.0     0.2803      0.223      1.254      0.210      -0.158       0.718"
Replace-InScope $table2 "year_2020.0    -0.1908      0.316     -0.605      0.546      -0.809       0.428" "year_2020.0    -0.5253      0.222     -2.369      0.018      -0.960      -0.091"
Replace-InScope $table2 "Omnibus:                     8518.972   Durbin-Watson:                   1.298" "Omnibus:                    96682.371   Durbin-Watson:                   1.314"
Replace-InScope $table2 "Prob(Omnibus):                  0.000   Jarque-Bera (JB):         45338131.041" "Prob(Omnibus):                  0.000   Jarque-Bera (JB):       8943850897.107"
Replace-InScope $table2 "Skew:                         -10.115   Prob(JB):                         0.00" "Skew:                          25.317   Prob(JB):                         0.00"
Replace-InScope $table2 "Kurtosis:                     461.676   Cond. No.                     1.25e+16" "Kurtosis:                    2337.738   Cond. No.                     1.17e+16"
Replace-InScope $table2 "[2] The smallest eigenvalue is 4.33e-27. This might indicate that there are" "[2] The smallest eigenvalue is 6.48e-25. This might indicate that there are"

# ---- Table 2 (epspx) new year rows, inserted after year_2020.0 row ----
InsertAfter-InScope $table2 "year_2020.0    -0.5253      0.222     -2.369      0.018      -0.960      -0.091" @(
    "year_2021.0     0.9938      0.217      4.584      0.000       0.569       1.419",
    "year_2022.0     0.1827      0.217      0.841      0.400      -0.243       0.608",
    "year_2023.0     0.6568      0.224      2.935      0.003       0.218       1.096"
)

Write-Host "Done."